$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "50.078.92"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +3.85%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.657.53"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +6.28%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "114.63"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +8.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "326.62"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.71%  "
$ws.Range("E7").Value = "  +2.02%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.557"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +3.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.30"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +5.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.12"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.75%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0826"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.76%  "
$ws.Range("E13").Value = "  +0.45%  "
$ws.Range("E14").Value = "  +4.00%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.073.26"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +6.23%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.653.06"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +5.51%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "50.010.27"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +3.94%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.28"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +3.49%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.78"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +2.92%  "
$ws.Range("E22").Value = "  +3.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "72.53"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.95%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "276.90"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +3.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.60"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +3.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.98"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +4.49%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.04"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +2.86%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "36.93"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +6.40%  "
$ws.Range("E30").Value = "  +1.52%  "
$ws.Range("E31").Value = "  +1.78%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "50.18"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.67%  "
$ws.Range("E33").Value = "  +3.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.69"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +2.83%  "
$ws.Range("E35").Value = "  +5.62%  "
$ws.Range("E36").Value = "  -0.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.01"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +9.08%  "
$ws.Range("E38").Value = "  +6.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.14"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +8.68%  "
$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.113"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.69%  "
$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "124.42"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "22.24"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.18%  "
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0319"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +5.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.105.44"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +5.09%  "
$ws.Range("E46").Value = "  +5.72%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.27"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +13.41%  "
$ws.Range("E48").Value = "  +3.87%  "
$ws.Range("E49").Value = "  +2.09%  "
$ws.Range("E50").Value = "  +3.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "60.38"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +6.33%  "
